$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value = "2025-03-26T12:33"
$ws.Range("B2").Value = "SHIFT_1"
$ws.Range("C2").Value = 200
$ws.Range("N2").Value = "Suriya"

# Row 3 updates
$ws.Range("A3").Value = "2025-03-27T12:22"
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = "Suri"
$ws.Range("O3").Value = "Suri"
$ws.Range("P3").Value = "Suri"

# Row 4 updates
$ws.Range("A4").Value = "2025-03-27T12:31"
$ws.Range("C4").Value = 100
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1
$ws.Range("N4").Value = "t"
$ws.Range("O4").Value = "t"
$ws.Range("P4").Value = "t"

# Row 5 new row - stored as inline/text strings for numeric-looking cells
$ws.Range("C5:M5").NumberFormat = "@"
$ws.Range("A5").Value = "2025-03-27T12:36"
$ws.Range("B5").Value = "SHIFT_1"
$ws.Range("C5").Value = "100"
$ws.Range("D5").Value = "1"
$ws.Range("E5").Value = "0"
$ws.Range("F5").Value = "0"
$ws.Range("G5").Value = "0"
$ws.Range("H5").Value = "0"
$ws.Range("I5").Value = "0"
$ws.Range("J5").Value = "0"
$ws.Range("K5").Value = "0"
$ws.Range("L5").Value = "0"
$ws.Range("M5").Value = "0"
$ws.Range("N5").Value = "t"
$ws.Range("O5").Value = "t"
$ws.Range("P5").Value = "t"
